$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo: "timetable" -> "timetables" (plural) in the single-user API URL
$ws.Range("A4").Value = "http://signmeinwebapi.azurewebsites.net/api/timetables/{id}"

# Update the selection to reflect where the author left the cursor after the fix
$ws.Range("A4").Select()
